$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.005514214498158991
$ws.Range("C2").Value = 0.007694010848738613
$ws.Range("D2").Value = 0.009578166933231158
$ws.Range("E2").Value = 0.01126444240632555

$ws.Range("B3").Value = 1.955155689973489
$ws.Range("C3").Value = 3.394831111512401
$ws.Range("D3").Value = 4.926115517268729
$ws.Range("E3").Value = 4.734337573985806

$ws.Range("B4").Value = -0.02220448350877621
$ws.Range("C4").Value = -0.02481912364874573
$ws.Range("D4").Value = -0.0324023320837963
$ws.Range("E4").ClearContents()

$ws.Range("B5").Value = -7.013272318267988
$ws.Range("C5").Value = -8.404164939436404
$ws.Range("D5").Value = -6.333969479823645
$ws.Range("E5").ClearContents()

$ws.Range("B6").Value = 0.01224697559470602
$ws.Range("C6").Value = 0.008368220489540192
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

$ws.Range("B7").Value = 2.795967234906044
$ws.Range("C7").Value = 2.273219974189968
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()

$ws.Range("B8").Value = 0.007657089661325975
$ws.Range("C8").Value = 0.008903211863784781
$ws.Range("D8").Value = 0.01063790674073911
$ws.Range("E8").Value = 0.01238057532362279

$ws.Range("B9").Value = 2.765221197521946
$ws.Range("C9").Value = 4.040882854193589
$ws.Range("D9").Value = 4.8014351120712
$ws.Range("E9").Value = 4.952453258095796

$ws.Range("B10").Value = -0.02017519830935221
$ws.Range("C10").Value = -0.0226542295666625
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()

$ws.Range("B11").Value = -6.339797393399877
$ws.Range("C11").Value = -7.389830631458113
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()

$ws.Range("B12").Value = 0.00835776495227241
$ws.Range("C12").Value = 0.004504963447397542
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()

$ws.Range("B13").Value = 1.895934115745534
$ws.Range("C13").Value = 1.119383893630111
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = 0.007078685298556125
$ws.Range("C14").Value = 0.008680487604028096
$ws.Range("D14").Value = 0.01065621721080113
$ws.Range("E14").Value = 0.01259666415202041

$ws.Range("B15").Value = 2.515867263755929
$ws.Range("C15").Value = 3.804121870816258
$ws.Range("D15").Value = 4.619232663797836
$ws.Range("E15").Value = 5.124331086033146

$ws.Range("B16").Value = -0.01773686138069842
$ws.Range("C16").Value = -0.02116683297095204
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()

$ws.Range("B17").Value = -5.513178630950585
$ws.Range("C17").Value = -6.616497594119018
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()

$ws.Range("B18").Value = 0.006810583529634118
$ws.Range("C18").Value = 0.004615601784072395
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()

$ws.Range("B19").Value = 1.622825916512491
$ws.Range("C19").Value = 1.349078961405825
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()

$ws.Range("B20").Value = 0.008852954431593705
$ws.Range("C20").Value = 0.01150173881096942
$ws.Range("D20").Value = 0.01358360093365739
$ws.Range("E20").Value = 0.01447989462478875

$ws.Range("B21").Value = 3.123567132121082
$ws.Range("C21").Value = 4.75913123278851
$ws.Range("D21").Value = 5.571118917644933
$ws.Range("E21").Value = 5.632052754125775

$ws.Range("B22").Value = -0.01652269895533423
$ws.Range("C22").Value = -0.01976221467413408
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()

$ws.Range("B23").Value = -4.936235049943778
$ws.Range("C23").Value = -5.837038295467231
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()

$ws.Range("B24").Value = 0.003779808589563728
$ws.Range("C24").Value = -0.0002637754810975874
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()

$ws.Range("B25").Value = 0.8851654475796069
$ws.Range("C25").Value = -0.07193459449966344
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
